$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 916.974
$ws.Cells.Item(1, 2).Value = 1023.369
$ws.Cells.Item(2, 1).Value = 934.543
$ws.Cells.Item(2, 2).Value = 1780.392
$ws.Cells.Item(2, 4).Value = -16.08
$ws.Cells.Item(3, 1).Value = 929.96
$ws.Cells.Item(3, 2).Value = 1675.658
$ws.Cells.Item(3, 4).Value = -13.95
$ws.Cells.Item(4, 1).Value = 926.124
$ws.Cells.Item(4, 2).Value = 1578.674
$ws.Cells.Item(4, 4).Value = -11.95
$ws.Cells.Item(5, 1).Value = 923.225
$ws.Cells.Item(5, 2).Value = 1482.496
$ws.Cells.Item(5, 4).Value = -9.93
$ws.Cells.Item(6, 1).Value = 920.535
$ws.Cells.Item(6, 2).Value = 1390.102
$ws.Cells.Item(6, 4).Value = -7.93
$ws.Cells.Item(7, 1).Value = 918.765
$ws.Cells.Item(7, 2).Value = 1297.492
$ws.Cells.Item(7, 4).Value = -5.95
$ws.Cells.Item(8, 1).Value = 917.58
$ws.Cells.Item(8, 2).Value = 1203.935
$ws.Cells.Item(8, 4).Value = -3.93
$ws.Cells.Item(9, 1).Value = 917.02
$ws.Cells.Item(9, 2).Value = 1111.692
$ws.Cells.Item(9, 4).Value = -1.93
$ws.Cells.Item(10, 1).Value = 916.899
$ws.Cells.Item(10, 2).Value = 1020.639
$ws.Cells.Item(10, 4).Value = 0.05
$ws.Cells.Item(11, 1).Value = 917.53
$ws.Cells.Item(11, 2).Value = 927.796
$ws.Cells.Item(11, 4).Value = 2.07
$ws.Cells.Item(12, 1).Value = 918.529
$ws.Cells.Item(12, 2).Value = 836.504
$ws.Cells.Item(12, 4).Value = 4.07
$ws.Cells.Item(13, 1).Value = 920.487
$ws.Cells.Item(13, 2).Value = 744.532
$ws.Cells.Item(13, 4).Value = 6.07
$ws.Cells.Item(14, 1).Value = 922.75
$ws.Cells.Item(14, 2).Value = 651.29
$ws.Cells.Item(14, 4).Value = 8.07
$ws.Cells.Item(15, 1).Value = 925.739
$ws.Cells.Item(15, 2).Value = 557.362
$ws.Cells.Item(15, 4).Value = 10.07
$ws.Cells.Item(16, 1).Value = 929.5
$ws.Cells.Item(16, 2).Value = 463.384
$ws.Cells.Item(16, 4).Value = 12.1
$ws.Cells.Item(17, 1).Value = 933.584
$ws.Cells.Item(17, 2).Value = 369.453
$ws.Cells.Item(17, 4).Value = 14.07
$ws.Cells.Item(18, 1).Value = 938.57
$ws.Cells.Item(18, 2).Value = 273.115
$ws.Cells.Item(18, 4).Value = 16.1
$ws.Cells.Item(19, 1).Value = 898.001
$ws.Cells.Item(19, 2).Value = 1023.342
$ws.Cells.Item(20, 1).Value = 915.276
$ws.Cells.Item(20, 2).Value = 1779.363
$ws.Cells.Item(20, 4).Value = -16.08
$ws.Cells.Item(21, 1).Value = 910.638
$ws.Cells.Item(21, 2).Value = 1674.001
$ws.Cells.Item(21, 4).Value = -13.95
$ws.Cells.Item(22, 1).Value = 907.042
$ws.Cells.Item(22, 2).Value = 1577.708
$ws.Cells.Item(22, 4).Value = -11.93
$ws.Cells.Item(23, 1).Value = 904.283
$ws.Cells.Item(23, 2).Value = 1483.639
$ws.Cells.Item(23, 4).Value = -9.95
$ws.Cells.Item(24, 1).Value = 901.598
$ws.Cells.Item(24, 2).Value = 1389.718
$ws.Cells.Item(24, 4).Value = -7.93
$ws.Cells.Item(25, 1).Value = 899.96
$ws.Cells.Item(25, 2).Value = 1296.457
$ws.Cells.Item(25, 4).Value = -5.95
$ws.Cells.Item(26, 1).Value = 898.671
$ws.Cells.Item(26, 2).Value = 1204.301
$ws.Cells.Item(26, 4).Value = -3.93
$ws.Cells.Item(27, 1).Value = 898.021
$ws.Cells.Item(27, 2).Value = 1112.383
$ws.Cells.Item(27, 4).Value = -1.93
$ws.Cells.Item(28, 1).Value = 898.001
$ws.Cells.Item(28, 2).Value = 1020.277
$ws.Cells.Item(28, 4).Value = 0.05
$ws.Cells.Item(29, 1).Value = 898.584
$ws.Cells.Item(29, 2).Value = 927.72
$ws.Cells.Item(29, 4).Value = 2.07
$ws.Cells.Item(30, 1).Value = 899.606
$ws.Cells.Item(30, 2).Value = 836.816
$ws.Cells.Item(30, 4).Value = 4.07
$ws.Cells.Item(31, 1).Value = 901.497
$ws.Cells.Item(31, 2).Value = 744.834
$ws.Cells.Item(31, 4).Value = 6.07
$ws.Cells.Item(32, 1).Value = 903.66
$ws.Cells.Item(32, 2).Value = 652.205
$ws.Cells.Item(32, 4).Value = 8.07
$ws.Cells.Item(33, 1).Value = 906.547
$ws.Cells.Item(33, 2).Value = 557.643
$ws.Cells.Item(33, 4).Value = 10.07
$ws.Cells.Item(34, 1).Value = 910.285
$ws.Cells.Item(34, 2).Value = 464.608
$ws.Cells.Item(34, 4).Value = 12.07
$ws.Cells.Item(35, 1).Value = 914.473
$ws.Cells.Item(35, 2).Value = 368.675
$ws.Cells.Item(35, 4).Value = 14.1
$ws.Cells.Item(36, 1).Value = 918.876
$ws.Cells.Item(36, 2).Value = 273.553
$ws.Cells.Item(36, 4).Value = 16.1
$ws.Cells.Item(37, 1).Value = 880.618
$ws.Cells.Item(37, 2).Value = 1023.821
$ws.Cells.Item(38, 1).Value = 897.167
$ws.Cells.Item(38, 2).Value = 1780.969
$ws.Cells.Item(38, 4).Value = -16.1
$ws.Cells.Item(39, 1).Value = 892.833
$ws.Cells.Item(39, 2).Value = 1675.692
$ws.Cells.Item(39, 4).Value = -13.97
$ws.Cells.Item(40, 1).Value = 889.273
$ws.Cells.Item(40, 2).Value = 1578.149
$ws.Cells.Item(40, 4).Value = -11.92
$ws.Cells.Item(41, 1).Value = 886.499
$ws.Cells.Item(41, 2).Value = 1483.727
$ws.Cells.Item(41, 4).Value = -9.92
$ws.Cells.Item(42, 1).Value = 884.193
$ws.Cells.Item(42, 2).Value = 1391.367
$ws.Cells.Item(42, 4).Value = -7.95
$ws.Cells.Item(43, 1).Value = 882.437
$ws.Cells.Item(43, 2).Value = 1297.131
$ws.Cells.Item(43, 4).Value = -5.92
$ws.Cells.Item(44, 1).Value = 881.213
$ws.Cells.Item(44, 2).Value = 1205.407
$ws.Cells.Item(44, 4).Value = -3.95
$ws.Cells.Item(45, 1).Value = 880.794
$ws.Cells.Item(45, 2).Value = 1112.39
$ws.Cells.Item(45, 4).Value = -1.92
$ws.Cells.Item(46, 1).Value = 880.623
$ws.Cells.Item(46, 2).Value = 1021.363
$ws.Cells.Item(46, 4).Value = 0.05
$ws.Cells.Item(47, 1).Value = 881.02
$ws.Cells.Item(47, 2).Value = 929.12
$ws.Cells.Item(47, 4).Value = 2.08
$ws.Cells.Item(48, 1).Value = 882.293
$ws.Cells.Item(48, 2).Value = 838.021
$ws.Cells.Item(48, 4).Value = 4.08
$ws.Cells.Item(49, 1).Value = 883.738
$ws.Cells.Item(49, 2).Value = 746.382
$ws.Cells.Item(49, 4).Value = 6.08
$ws.Cells.Item(50, 1).Value = 886.017
$ws.Cells.Item(50, 2).Value = 654.025
$ws.Cells.Item(50, 4).Value = 8.08
$ws.Cells.Item(51, 1).Value = 888.948
$ws.Cells.Item(51, 2).Value = 560.247
$ws.Cells.Item(51, 4).Value = 10.05
$ws.Cells.Item(52, 1).Value = 892.294
$ws.Cells.Item(52, 2).Value = 466.099
$ws.Cells.Item(52, 4).Value = 12.1
$ws.Cells.Item(53, 1).Value = 896.428
$ws.Cells.Item(53, 2).Value = 370.87
$ws.Cells.Item(53, 4).Value = 14.1
$ws.Cells.Item(54, 1).Value = 900.605
$ws.Cells.Item(54, 2).Value = 275.111
$ws.Cells.Item(54, 4).Value = 16.08
$ws.Cells.Item(55, 1).Value = 864.266
$ws.Cells.Item(55, 2).Value = 1023.289
$ws.Cells.Item(56, 1).Value = 880.319
$ws.Cells.Item(56, 2).Value = 1778.195
$ws.Cells.Item(56, 4).Value = -16.08
$ws.Cells.Item(57, 1).Value = 876.005
$ws.Cells.Item(57, 2).Value = 1673.183
$ws.Cells.Item(57, 4).Value = -13.98
$ws.Cells.Item(58, 1).Value = 872.679
$ws.Cells.Item(58, 2).Value = 1575.632
$ws.Cells.Item(58, 4).Value = -11.93
$ws.Cells.Item(59, 1).Value = 870.102
$ws.Cells.Item(59, 2).Value = 1482.053
$ws.Cells.Item(59, 4).Value = -9.93
$ws.Cells.Item(60, 1).Value = 867.614
$ws.Cells.Item(60, 2).Value = 1389.088
$ws.Cells.Item(60, 4).Value = -7.95
$ws.Cells.Item(61, 1).Value = 866.039
$ws.Cells.Item(61, 2).Value = 1295.472
$ws.Cells.Item(61, 4).Value = -5.93
$ws.Cells.Item(62, 1).Value = 865.02
$ws.Cells.Item(62, 2).Value = 1204.119
$ws.Cells.Item(62, 4).Value = -3.93
$ws.Cells.Item(63, 1).Value = 864.416
$ws.Cells.Item(63, 2).Value = 1112.955
$ws.Cells.Item(63, 4).Value = -1.95
$ws.Cells.Item(64, 1).Value = 864.173
$ws.Cells.Item(64, 2).Value = 1019.616
$ws.Cells.Item(64, 4).Value = 0.07
$ws.Cells.Item(65, 1).Value = 865
$ws.Cells.Item(65, 2).Value = 929.235
$ws.Cells.Item(65, 4).Value = 2.07
$ws.Cells.Item(66, 1).Value = 865.987
$ws.Cells.Item(66, 2).Value = 837.896
$ws.Cells.Item(66, 4).Value = 4.07
$ws.Cells.Item(67, 1).Value = 867.521
$ws.Cells.Item(67, 2).Value = 746.024
$ws.Cells.Item(67, 4).Value = 6.07
$ws.Cells.Item(68, 1).Value = 869.602
$ws.Cells.Item(68, 2).Value = 653.378
$ws.Cells.Item(68, 4).Value = 8.07
$ws.Cells.Item(69, 1).Value = 872.358
$ws.Cells.Item(69, 2).Value = 559.993
$ws.Cells.Item(69, 4).Value = 10.07
$ws.Cells.Item(70, 1).Value = 875.925
$ws.Cells.Item(70, 2).Value = 466.413
$ws.Cells.Item(70, 4).Value = 12.07
$ws.Cells.Item(71, 1).Value = 879.586
$ws.Cells.Item(71, 2).Value = 370.722
$ws.Cells.Item(71, 4).Value = 14.1
$ws.Cells.Item(72, 1).Value = 883.928
$ws.Cells.Item(72, 2).Value = 275.893
$ws.Cells.Item(72, 4).Value = 16.1
$ws.Cells.Item(73, 1).Value = 849.223
$ws.Cells.Item(73, 2).Value = 1023.567
$ws.Cells.Item(74, 1).Value = 864.551
$ws.Cells.Item(74, 2).Value = 1777.473
$ws.Cells.Item(74, 4).Value = -16.08
$ws.Cells.Item(75, 1).Value = 860.626
$ws.Cells.Item(75, 2).Value = 1672.561
$ws.Cells.Item(75, 4).Value = -13.98
$ws.Cells.Item(76, 1).Value = 857.431
$ws.Cells.Item(76, 2).Value = 1575.263
$ws.Cells.Item(76, 4).Value = -11.93
$ws.Cells.Item(77, 1).Value = 854.68
$ws.Cells.Item(77, 2).Value = 1481.736
$ws.Cells.Item(77, 4).Value = -9.93
$ws.Cells.Item(78, 1).Value = 852.561
$ws.Cells.Item(78, 2).Value = 1388.454
$ws.Cells.Item(78, 4).Value = -7.93
$ws.Cells.Item(79, 1).Value = 851
$ws.Cells.Item(79, 2).Value = 1295.755
$ws.Cells.Item(79, 4).Value = -5.95
$ws.Cells.Item(80, 1).Value = 850
$ws.Cells.Item(80, 2).Value = 1203.581
$ws.Cells.Item(80, 4).Value = -3.93
$ws.Cells.Item(81, 1).Value = 849.15
$ws.Cells.Item(81, 2).Value = 1112.662
$ws.Cells.Item(81, 4).Value = -1.93
$ws.Cells.Item(82, 1).Value = 849.222
$ws.Cells.Item(82, 2).Value = 1020.99
$ws.Cells.Item(82, 4).Value = 0.05
$ws.Cells.Item(83, 1).Value = 849.994
$ws.Cells.Item(83, 2).Value = 929.325
$ws.Cells.Item(83, 4).Value = 2.07
$ws.Cells.Item(84, 1).Value = 850.987
$ws.Cells.Item(84, 2).Value = 838.124
$ws.Cells.Item(84, 4).Value = 4.07
$ws.Cells.Item(85, 1).Value = 852.484
$ws.Cells.Item(85, 2).Value = 746.507
$ws.Cells.Item(85, 4).Value = 6.07
$ws.Cells.Item(86, 1).Value = 854.503
$ws.Cells.Item(86, 2).Value = 653.227
$ws.Cells.Item(86, 4).Value = 8.07
$ws.Cells.Item(87, 1).Value = 857.047
$ws.Cells.Item(87, 2).Value = 560.664
$ws.Cells.Item(87, 4).Value = 10.07
$ws.Cells.Item(88, 1).Value = 860.429
$ws.Cells.Item(88, 2).Value = 467.003
$ws.Cells.Item(88, 4).Value = 12.07
$ws.Cells.Item(89, 1).Value = 864.305
$ws.Cells.Item(89, 2).Value = 372.666
$ws.Cells.Item(89, 4).Value = 14.1
$ws.Cells.Item(90, 1).Value = 868.406
$ws.Cells.Item(90, 2).Value = 277.837
$ws.Cells.Item(90, 4).Value = 16.07

# Clear rows 91-108 (the last C=2999.3 block) entirely
$ws.Range("A91:D108").ClearContents()

# Update selection to match new range A1:D90
$ws.Range("A1:D90").Select()